# Update FlashScore odds/stats figures for 2024-10-12 workbook.
# Each assignment below corresponds to one changed value from the
# authoritative diff between the previous and current OOXML.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 1.61

# Row 6
$ws.Range("G6").Value = 2.4
$ws.Range("R6").Value = 1.5

# Row 8
$ws.Range("O8").Value = 1.25
$ws.Range("P8").Value = 3.75
$ws.Range("Q8").Value = 1.9
$ws.Range("R8").Value = 1.95

# Row 13
$ws.Range("N13").Value = 17

# Row 19
$ws.Range("K19").Value = 1.87

# Row 21
$ws.Range("J21").Value = 2.37
$ws.Range("Q21").Value = 2.05
$ws.Range("R21").Value = 1.75
$ws.Range("W21").Value = 6.5
$ws.Range("AD21").Value = 6.5
$ws.Range("AE21").Value = 17
$ws.Range("AF21").Value = 51
$ws.Range("AH21").Value = 23
$ws.Range("AL21").Value = 41
$ws.Range("AX21").Value = 26
$ws.Range("AY21").Value = 34

# Row 22
$ws.Range("R22").Value = 1.63

# Row 23
$ws.Range("I23").Value = 2.87
$ws.Range("R23").Value = 1.62

# Row 24
$ws.Range("R24").Value = 1.72

# Row 25
$ws.Range("G25").Value = 2.3
$ws.Range("H25").Value = 3.1
$ws.Range("I25").Value = 3.4
$ws.Range("J25").Value = 3
$ws.Range("L25").Value = 4
$ws.Range("R25").Value = 1.62
$ws.Range("X25").Value = 10
$ws.Range("Y25").Value = 9.5
$ws.Range("AC25").Value = 8
$ws.Range("AV25").Value = 51
$ws.Range("AZ25").Value = 67

# Row 26
$ws.Range("G26").Value = 1.38
$ws.Range("BD26").Value = 151

# Row 35
$ws.Range("G35").Value = 1.82
$ws.Range("L35").Value = 4.45
$ws.Range("O35").Value = 1.3
$ws.Range("P35").Value = 2.92
$ws.Range("Q35").Value = 1.95
$ws.Range("V35").Value = 1.83
$ws.Range("AF35").Value = 70
$ws.Range("AG35").Value = 11.75
$ws.Range("AH35").Value = 25
$ws.Range("AJ35").Value = 75
$ws.Range("AU35").Value = 6.9
$ws.Range("AY35").Value = 27
$ws.Range("BB35").Value = 300

# Row 36
$ws.Range("G36").Value = 3
$ws.Range("H36").Value = 3.1
$ws.Range("J36").Value = 3.5
$ws.Range("K36").Value = 2.02
$ws.Range("L36").Value = 2.95
$ws.Range("O36").Value = 1.28
$ws.Range("P36").Value = 3.05
$ws.Range("S36").Value = 1.4
$ws.Range("T36").Value = 2.52
$ws.Range("W36").Value = 10
$ws.Range("X36").Value = 17
$ws.Range("Z36").Value = 40
$ws.Range("AB36").Value = 29
$ws.Range("AD36").Value = 6
$ws.Range("AE36").Value = 12
$ws.Range("AF36").Value = 50
$ws.Range("AG36").Value = 8
$ws.Range("AH36").Value = 11.75
$ws.Range("AI36").Value = 9
$ws.Range("AJ36").Value = 25
$ws.Range("AK36").Value = 18.5
$ws.Range("AL36").Value = 26
$ws.Range("AM36").Value = 350
$ws.Range("AP36").Value = 22
$ws.Range("AQ36").Value = 75
$ws.Range("AR36").Value = 100
$ws.Range("AS36").Value = 250
$ws.Range("AT36").Value = 2.5
$ws.Range("AX36").Value = 12.5
$ws.Range("AY36").Value = 20
$ws.Range("BA36").Value = 80

# Row 37
$ws.Range("H37").Value = 3.25
$ws.Range("I37").Value = 3
$ws.Range("K37").Value = 2.05
$ws.Range("L37").Value = 3.55
$ws.Range("P37").Value = 2.62
$ws.Range("U37").Value = 1.87
$ws.Range("W37").Value = 6.7
$ws.Range("AA37").Value = 20
$ws.Range("AD37").Value = 6.3
$ws.Range("AE37").Value = 16.5
$ws.Range("AG37").Value = 8
$ws.Range("AT37").Value = 2.42
$ws.Range("AW37").Value = 4.75
$ws.Range("AX37").Value = 16.5
$ws.Range("BB37").Value = 350
